$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E7").Value = 16.831
$ws.Range("B9").Value = 6.345999999999999
$ws.Range("E12").Value = 17.646
$ws.Range("E14").Value = 16.866
$ws.Range("B18").Value = 5.029999999999999
$ws.Range("B20").Value = 6.239999999999999
$ws.Range("E26").Value = 16.525
$ws.Range("B27").Value = 5.752
$ws.Range("E27").Value = 16.672
$ws.Range("E29").Value = 17
$ws.Range("B35").Value = 7.657999999999999
$ws.Range("E37").Value = 16.855
$ws.Range("E38").Value = 16.871
$ws.Range("E51").Value = 16.65
$ws.Range("E52").Value = 16.543
$ws.Range("E55").Value = 16.494
$ws.Range("B69").Value = 5.627000000000001
$ws.Range("E69").Value = 17.321
$ws.Range("E70").Value = 17.559
$ws.Range("B76").Value = 6.308
$ws.Range("B78").Value = 8.254
$ws.Range("E81").Value = 16.457
$ws.Range("B82").Value = 5.366
$ws.Range("B83").Value = 5.88
$ws.Range("E83").Value = 16.917
$ws.Range("B93").Value = 5.724
$ws.Range("E102").Value = 16.724
